# Weekly update: insert 4 new measurement rows at the top of the
# "Vega Monumental Concepción - Repollo" block (row 136) and push the
# existing rows down, so the newest week's entries lead the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 136; everything from 136..202
# shifts down to 140..206 (dimension grows from R202 to R206).
$ws.Rows.Item(136).Resize(4).Insert()

# Fill in the 4 newly-inserted rows with this week's data.
$newRows = @(
    @{ Row = 136; H = "Crespo record"; I = "Primera"; J = 2000; K = 700; L = 800; M = 750; O = "Región Metropolitana"; P = 750 },
    @{ Row = 137; H = "Crespo record"; I = "Segunda"; J = 1000; K = 600; L = 600; M = 600; O = "Región Metropolitana"; P = 600 },
    @{ Row = 138; H = "Morada(o)";     I = "Primera"; J = 600;  K = 700; L = 800; M = 750; O = "Región Metropolitana"; P = 750 },
    @{ Row = 139; H = "Morada(o)";     I = "Segunda"; J = 300;  K = 600; L = 600; M = 600; O = "Región Metropolitana"; P = 600 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = 44468
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = 100112006
    $ws.Cells.Item($row, 7).Value = "Repollo"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
